$d = $word.ActiveDocument

$replacements = @(
    @("2025-01-29 Wednesday", "2025-01-30 Thursday"),
    @("76×22=", "23×84="),
    @("99×37=", "62×82="),
    @("27×30=", "85×13="),
    @("95×78=", "71×88="),
    @("62×86=", "47×81="),
    @("83×32=", "35×56="),
    @("11×72=", "25×61="),
    @("79×26=", "52×44="),
    @("68×86=", "19×19="),
    @("16×28=", "83×74="),
    @("59×80=", "81×49="),
    @("41×89=", "34×66="),
    @("73×62=", "97×58="),
    @("34×64=", "33×62="),
    @("38×56=", "20×77="),
    @("54×63=", "73×41="),
    @("35×61=", "63×64="),
    @("33×30=", "92×20="),
    @("77×68=", "15×55="),
    @("95×36=", "73×51="),
    @("25×66=", "43×47="),
    @("72×45=", "69×16="),
    @("42×87=", "46×26="),
    @("15×49=", "26×45="),
    @("93×44=", "61×86=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
